# Auto-generated Excel COM-interop script applying the diff changes
# (commit: Update latest output (run 144))

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update rows 3 and 4 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("A3").Value = 46052.04166666666
$wsSchedule.Range("B3").Value = 46052.20833333334
$wsSchedule.Range("E3").Value = 551.7321029999999
$wsSchedule.Range("F3").Value = 36.49021845238095
$wsSchedule.Range("A4").Value = 46052.29166666666
$wsSchedule.Range("B4").Value = 46052.625
$wsSchedule.Range("E4").Value = 641.5914375
$wsSchedule.Range("F4").Value = 21.21664806547619

# --- Sheet "Detailed": update rows 38-97 ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B38").Value = 57.06007
$wsDetailed.Range("B39").Value = 73.19
$wsDetailed.Range("B40").Value = 92.30755000000001
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 96.67016
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 101.25
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 90.83913
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 78
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 73.19
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 106.9792
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 105
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 92.96397
$wsDetailed.Range("B50").Value = 93.8609
$wsDetailed.Range("E50").Value = "OFF"
$wsDetailed.Range("B51").Value = 78
$wsDetailed.Range("E51").Value = "OFF"
$wsDetailed.Range("B52").Value = 76.58247
$wsDetailed.Range("B53").Value = 64.89
$wsDetailed.Range("B54").Value = 66.90913
$wsDetailed.Range("B55").Value = 66.60672
$wsDetailed.Range("B56").Value = 68.34329
$wsDetailed.Range("B57").Value = 71.34737
$wsDetailed.Range("B58").Value = 73.20010000000001
$wsDetailed.Range("E58").Value = "ON"
$wsDetailed.Range("B59").Value = 78
$wsDetailed.Range("E59").Value = "ON"
$wsDetailed.Range("B60").Value = 86.56122999999999
$wsDetailed.Range("B61").Value = 108.89
$wsDetailed.Range("B62").Value = 127.50518
$wsDetailed.Range("B63").Value = 108.89
$wsDetailed.Range("B64").Value = 68.43391
$wsDetailed.Range("E64").Value = "ON"
$wsDetailed.Range("B65").Value = 56.98
$wsDetailed.Range("B66").Value = 36.11552
$wsDetailed.Range("B67").Value = 49.44934
$wsDetailed.Range("B68").Value = 33.9499
$wsDetailed.Range("B70").Value = 36.0601
$wsDetailed.Range("B71").Value = 31.49306
$wsDetailed.Range("B73").Value = 36.07
$wsDetailed.Range("B75").Value = 36.0601
$wsDetailed.Range("B78").Value = 36.07
$wsDetailed.Range("B79").Value = 57.06007
$wsDetailed.Range("B80").Value = 73.05861
$wsDetailed.Range("E80").Value = "OFF"
$wsDetailed.Range("B81").Value = 73.20007
$wsDetailed.Range("B82").Value = 84.79000000000001
$wsDetailed.Range("B83").Value = 68.44540000000001
$wsDetailed.Range("B84").Value = 73.88037
$wsDetailed.Range("B85").Value = 101.25
$wsDetailed.Range("B86").Value = 108.89
$wsDetailed.Range("B87").Value = 125.64596
$wsDetailed.Range("B88").Value = 283.96
$wsDetailed.Range("B91").Value = 299.99
$wsDetailed.Range("B92").Value = 240.89
$wsDetailed.Range("B93").Value = 147.89
$wsDetailed.Range("B94").Value = 138.21589
$wsDetailed.Range("B95").Value = 147.51
$wsDetailed.Range("B96").Value = 158.64449
$wsDetailed.Range("B97").Value = 133.70206
